# edit.ps1
# Applies a row-wise shuffle of the "Fecha, Calidad, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Unidad de comercializacion,
# Precio $/Kg, Kg / unidad" columns (D, L, M, N, O, P, Q, S, T) across data rows
# 2-20 of the active worksheet, per the commit "Fruta / hortaliza, semanal".
# Columns A,B,C,E,F,G,H,I,J,K,R are identical across all rows and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source snapshot of rows 2-20 for columns D,L,M,N,O,P,Q,S,T (from the original workbook)
$rowsData = @{}
$rowsData[2] = @{ D = 44307; L = 'Primera'; M = 250; N = 19000; O = 20000; P = 19500; Q = '$/bandeja 18 kilos'; S = 1083; T = 18 }
$rowsData[3] = @{ D = 44789; L = 'Segunda'; M = 250; N = 19000; O = 20000; P = 19500; Q = '$/bandeja 18 kilos'; S = 1083; T = 18 }
$rowsData[4] = @{ D = 44263; L = 'Primera'; M = 250; N = 21000; O = 22000; P = 21500; Q = '$/caja 18 kilos'; S = 1194; T = 18 }
$rowsData[5] = @{ D = 44491; L = 'Primera'; M = 300; N = 14000; O = 15000; P = 14500; Q = '$/bandeja 10 kilos'; S = 1450; T = 10 }
$rowsData[6] = @{ D = 44602; L = 'Primera'; M = 270; N = 20000; O = 21000; P = 20500; Q = '$/bandeja 18 kilos'; S = 1139; T = 18 }
$rowsData[7] = @{ D = 44673; L = 'Especial'; M = 400; N = 14000; O = 15000; P = 14500; Q = '$/bandeja 10 kilos'; S = 1450; T = 10 }
$rowsData[8] = @{ D = 44614; L = 'Primera'; M = 250; N = 20000; O = 21000; P = 20500; Q = '$/bandeja 18 kilos'; S = 1139; T = 18 }
$rowsData[9] = @{ D = 44991; L = 'Primera'; M = 250; N = 24000; O = 25000; P = 24500; Q = '$/bandeja 18 kilos'; S = 1361; T = 18 }
$rowsData[10] = @{ D = 44487; L = 'Primera'; M = 300; N = 14000; O = 15000; P = 14500; Q = '$/bandeja 10 kilos'; S = 1450; T = 10 }
$rowsData[11] = @{ D = 44629; L = 'Segunda'; M = 300; N = 17000; O = 18000; P = 17500; Q = '$/bandeja 18 kilos'; S = 972; T = 18 }
$rowsData[12] = @{ D = 44616; L = 'Segunda'; M = 300; N = 16000; O = 17000; P = 16500; Q = '$/caja 18 kilos granel'; S = 917; T = 18 }
$rowsData[13] = @{ D = 44656; L = 'Primera'; M = 270; N = 19000; O = 20000; P = 19500; Q = '$/bandeja 18 kilos'; S = 1083; T = 18 }
$rowsData[14] = @{ D = 44291; L = 'Primera'; M = 200; N = 17000; O = 18000; P = 17500; Q = '$/bandeja 18 kilos'; S = 972; T = 18 }
$rowsData[15] = @{ D = 44706; L = 'Primera'; M = 400; N = 9000; O = 10000; P = 9500; Q = '$/bandeja 10 kilos'; S = 950; T = 10 }
$rowsData[16] = @{ D = 44418; L = 'Primera'; M = 240; N = 10000; O = 11000; P = 10500; Q = '$/bandeja 10 kilos'; S = 1050; T = 10 }
$rowsData[17] = @{ D = 44489; L = 'Primera'; M = 300; N = 26000; O = 27000; P = 26500; Q = '$/bandeja 18 kilos'; S = 1472; T = 18 }
$rowsData[18] = @{ D = 44323; L = 'Primera'; M = 270; N = 21000; O = 22000; P = 21500; Q = '$/bandeja 18 kilos'; S = 1194; T = 18 }
$rowsData[19] = @{ D = 44784; L = 'Primera'; M = 300; N = 19000; O = 20000; P = 19500; Q = '$/bandeja 18 kilos'; S = 1083; T = 18 }
$rowsData[20] = @{ D = 44819; L = 'Primera'; M = 300; N = 17000; O = 18000; P = 17500; Q = '$/bandeja 10 kilos'; S = 1750; T = 10 }

# Mapping: new row -> source row (from original data) to pull values from
$order = @{}
$order[2] = 15
$order[3] = 13
$order[4] = 5
$order[5] = 11
$order[6] = 10
$order[7] = 18
$order[8] = 2
$order[9] = 9
$order[10] = 6
$order[11] = 16
$order[12] = 19
$order[13] = 12
$order[14] = 17
$order[15] = 4
$order[16] = 8
$order[17] = 20
$order[18] = 3
$order[19] = 14
$order[20] = 7

# Column letters to number mapping used below
$colNums = @{ D = 4; L = 12; M = 13; N = 14; O = 15; P = 16; Q = 17; S = 19; T = 20 }

foreach ($destRow in ($order.Keys | Sort-Object)) {
    $srcRow = $order[$destRow]
    $data = $rowsData[$srcRow]

    $ws.Cells.Item($destRow, $colNums.D).Value = $data.D
    $ws.Cells.Item($destRow, $colNums.L).Value = $data.L
    $ws.Cells.Item($destRow, $colNums.M).Value = $data.M
    $ws.Cells.Item($destRow, $colNums.N).Value = $data.N
    $ws.Cells.Item($destRow, $colNums.O).Value = $data.O
    $ws.Cells.Item($destRow, $colNums.P).Value = $data.P
    $ws.Cells.Item($destRow, $colNums.Q).Value = $data.Q
    $ws.Cells.Item($destRow, $colNums.S).Value = $data.S
    $ws.Cells.Item($destRow, $colNums.T).Value = $data.T
}
